$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: find the single paragraph whose text contains a given marker
# substring (and, optionally, does NOT contain an exclusion substring -
# needed to disambiguate the two "FOREIGN KEY (countryID) ..." lines,
# one of which ends in a comma and one of which doesn't).
# ---------------------------------------------------------------------
function Find-ParagraphContaining($needle, $excludeNeedle) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        if ($t -like "*$needle*") {
            if ($excludeNeedle -and ($t -like "*$excludeNeedle*")) {
                continue
            }
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Helper: force a clean run split at a zero-length Range position by
# adding (and immediately removing) a throwaway bookmark there. Adding
# a bookmark always breaks the enclosing run into two runs at that
# character offset without introducing any residual run formatting
# (no stray <w:rPr/>), which is exactly what's needed to reproduce the
# "typo fix split into separate runs" shape seen in the target diff.
# ---------------------------------------------------------------------
function Split-RunAt($range) {
    $d.Bookmarks.Add("zzTmpSplit", $range) | Out-Null
    $d.Bookmarks("zzTmpSplit").Delete()
}

# ---------------------------------------------------------------------
# 1) "FOREIGH KEY (locationID) REFERENCE location (locationID)"
#    -> fix the "FOREIGH" typo to "FOREIGN", ending up split across
#       three runs: "FOREIG" | "N" | " KEY (locationID) ...".
# ---------------------------------------------------------------------
$p1 = Find-ParagraphContaining "FOREIGH KEY (locationID)" $null
$p1Start = $p1.Range.Start
$hRange = $d.Range($p1Start + 6, $p1Start + 7)
$hRange.Text = "N"
Split-RunAt($d.Range($p1Start + 6, $p1Start + 6))
Split-RunAt($d.Range($p1Start + 7, $p1Start + 7))

# ---------------------------------------------------------------------
# 2) Collapse the "FOREIGN KEY (countryID) REFERENCE country
#    (countryID)," paragraph (currently 4 runs) back down into a
#    single run. A Find/Replace over the whole paragraph with
#    identical find/replacement text re-coalesces the runs without
#    altering the visible text. (Range.Text carries a trailing CR
#    paragraph mark that must be stripped before using it as literal
#    Find/Replacement text.)
# ---------------------------------------------------------------------
$p2 = Find-ParagraphContaining "country (countryID)," $null
$p2Text = $p2.Range.Text.TrimEnd([char]13)
$f2 = $p2.Range.Find
$f2.Text = $p2Text
$f2.Execute($null, $false, $false, $false, $false, $false, $true, 0, $false, $p2Text, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Same collapse for "FOREIGN KEY (ratingID) REFERENCE rating
#    (ratingID),".
# ---------------------------------------------------------------------
$p3 = Find-ParagraphContaining "rating (ratingID)," $null
$p3Text = $p3.Range.Text.TrimEnd([char]13)
$f3 = $p3.Range.Find
$f3.Text = $p3Text
$f3.Execute($null, $false, $false, $false, $false, $false, $true, 0, $false, $p3Text, 2) | Out-Null

# ---------------------------------------------------------------------
# 4) "FOREIGH KEY (itemID) REFERENCE item(itemID)"
#    -> fix "FOREIGH" to "FOREIGN", split into "FOREIG" | "N" |
#       " KEY (itemID) ...", and leave a "_GoBack" bookmark sitting
#       right after the "N" (between runs 2 and 3), matching the
#       last-edit-position bookmark Word drops at the point of the
#       most recent text change.
# ---------------------------------------------------------------------
$p4 = Find-ParagraphContaining "FOREIGH KEY (itemID)" $null
$p4Start = $p4.Range.Start
$hRange2 = $d.Range($p4Start + 6, $p4Start + 7)
$hRange2.Text = "N"
Split-RunAt($d.Range($p4Start + 6, $p4Start + 6))
$d.Bookmarks.Add("_GoBack", $d.Range($p4Start + 7, $p4Start + 7)) | Out-Null
